$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Labels for the expanded correlation matrix (columns B..I / rows 2..9)
$labels = @("S&P500", "Gold", "Kospi200", "USD", "WTI", "K_treasury", "K_corp_bond", "global_bonds")

# Full symmetric correlation matrix values (8x8), row i / col j aligning with $labels
$values = @(
    @(1, 0.6998666321284442, 0.8150783844262781, 0.1156948483490844, -0.6037594671456528, 0.784067051792052, 0.7495275017817324, 0.9202057654033214),
    @(0.6998666321284442, 1, 0.8997079696946566, 0.04379182799911028, -0.5314214205653955, 0.6190401330723398, 0.7153791319941667, 0.5879142527021959),
    @(0.8150783844262781, 0.8997079696946566, 1, -0.1459691394944984, -0.5005338126282869, 0.6707553186115066, 0.7155247127037971, 0.7020722691999959),
    @(0.1156948483490844, 0.04379182799911028, -0.1459691394944984, 1, -0.4834998051590432, 0.1001054865159668, 0.04754038645372555, 0.1860429897050268),
    @(-0.6037594671456528, -0.5314214205653955, -0.5005338126282869, -0.4834998051590432, 1, -0.7454485727909008, -0.6810985214977375, -0.5322923581610045),
    @(0.784067051792052, 0.6190401330723398, 0.6707553186115066, 0.1001054865159668, -0.7454485727909008, 1, 0.9182158444064913, 0.6773559548205431),
    @(0.7495275017817324, 0.7153791319941667, 0.7155247127037971, 0.04754038645372555, -0.6810985214977375, 0.9182158444064913, 1, 0.6073162434477504),
    @(0.9202057654033214, 0.5879142527021959, 0.7020722691999959, 0.1860429897050268, -0.5322923581610045, 0.6773559548205431, 0.6073162434477504, 1)
)

# Reference cells that already carry the bold/centered/bordered header style
$headerFormatSrc = $ws.Range("B1")
$labelFormatSrc = $ws.Range("A2")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $col = $i + 2   # B=2 .. I=9
    $row = $i + 2   # 2..9

    # Column header on row 1
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = $labels[$i]
    if ($col -gt 4) {
        $headerFormatSrc.Copy()
        $headerCell.PasteSpecial(-4122)
    }

    # Row label in column A
    $labelCell = $ws.Cells.Item($row, 1)
    $labelCell.Value = $labels[$i]
    if ($row -gt 4) {
        $labelFormatSrc.Copy()
        $labelCell.PasteSpecial(-4122)
    }

    # Correlation values for this row across columns B..I
    for ($j = 0; $j -lt $labels.Length; $j++) {
        $ws.Cells.Item($row, $j + 2).Value = $values[$i][$j]
    }
}

$excel.CutCopyMode = $false
